$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns with refreshed market data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.82%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.79%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.474"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.93%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08085"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.63%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.46%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.293"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.08%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.896"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.57%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9413"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.90%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1214"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.99%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1893"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.04%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09728"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.96%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04140"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6.60%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1068"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.70%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001270"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.55%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006053"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.90%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.571"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.28%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.521"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.92%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1347"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.14%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2494"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.61%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04374"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.87%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001236"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.86%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004292"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.52%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001236"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.80%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004002"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.26%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02655"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.93%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05454"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.16%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007640"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.90%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009724"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.21%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.86%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002122"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.44%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009902"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.13%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007116"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.63%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000754"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.41%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003555"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.50%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002277"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.08%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002110"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.41%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.41%"
